$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd negative sign on the end year for the Student Research Grant row
$ws.Range("D2").Value = 2021

# The end year for that same row is actually a continuation string ",2022"
# (not a standalone numeric year), so store it as text.
$ws.Range("E2").Value = ",2022"

# Reflect the last user selection having moved to E2
$ws.Range("E2").Select()
